# ControlLab/Lab4/raw_data.xlsx - "raw data to construct the bode plot"
#
# This script reproduces, via Excel COM automation, the edits captured in the
# target diff: a handful of corrected/added source values in columns B, C, D
# and G, a new block of Phase-2 formulas built from columns A-C mirrored into
# new columns M, N, O, and a refreshed set of "deltaT"/Phase formulas (col I)
# for the rows where a deltaT (col G) reading was filled in.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Corrections to existing measured data
# ---------------------------------------------------------------------------
$ws.Range("B6").Value2 = 3.91
$ws.Range("C17").Value2 = 5.2
$ws.Range("C32").Value2 = 1.8

# ---------------------------------------------------------------------------
# 2. Newly filled-in deltaT (G) readings for the higher frequency rows, and a
#    correction to the one that already existed (G30).
# ---------------------------------------------------------------------------
$ws.Range("G21").Value2 = 10.4
$ws.Range("G22").Value2 = 10.4
$ws.Range("G23").Value2 = 10
$ws.Range("G24").Value2 = 9.6
$ws.Range("G25").Value2 = 9.6
$ws.Range("G26").Value2 = 9.2
$ws.Range("G27").Value2 = 8.8
$ws.Range("G28").Value2 = 8.4
$ws.Range("G29").Value2 = 7.6
$ws.Range("G30").Value2 = 7.2

# ---------------------------------------------------------------------------
# 3. Newly filled-in Vout*2 (D) readings for rows 30-34.
# ---------------------------------------------------------------------------
$ws.Range("D30").Value2 = 116
$ws.Range("D31").Value2 = 118
$ws.Range("D32").Value2 = 128
$ws.Range("D33").Value2 = 140
$ws.Range("D34").Value2 = 144

# ---------------------------------------------------------------------------
# 4. Rows 21-29 now have a deltaT reading, so their Phase (I) formula switches
#    from the old "D*-1" style to the "-A*G*10^-3*360" style already used by
#    the rest of the table (rows 7-20, 30-37).
# ---------------------------------------------------------------------------
$ws.Range("I21").Formula = "=-A21*G21*10^-3*360"
$ws.Range("I22:I29").Formula = "=-A22*G22*10^-3*360"

# ---------------------------------------------------------------------------
# 5. New "Phase 2" helper columns: M (= freq, pasted as values), N (=
#    Vout/Vin, a live formula) and O (= Phase, pasted as values) for every
#    data row, 2 through 37.
# ---------------------------------------------------------------------------
$ws.Range("N2").Formula = "=C2/B2"
$ws.Range("N3:N37").Formula = "=C3/B3"

for ($r = 2; $r -le 37; $r++) {
    $ws.Range("M$r").Value2 = $ws.Range("A$r").Value2
    $ws.Range("O$r").Value2 = $ws.Range("I$r").Value2
}

# ---------------------------------------------------------------------------
# 6. Cursor/selection as left by the author.
# ---------------------------------------------------------------------------
$ws.Range("G31").Select()
